$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.311.08'
$ws.Range("E2").Value = '  -0.27%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.874.46'
$ws.Range("E3").Value = '  -0.13%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.18%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7093'
$ws.Range("E5").Value = '  -0.78%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.12'
$ws.Range("E6").Value = '  +0.01%  '

# Row 7
$ws.Range("E7").Value = '  +0.09%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07786'
$ws.Range("E8").Value = '  +0.38%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3106'
$ws.Range("E9").Value = '  -0.49%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '25.07'
$ws.Range("E10").Value = '  +0.50%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08412'
$ws.Range("E11").Value = '  +0.35%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.854.27'
$ws.Range("E12").Value = '  -2.09%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.242'
$ws.Range("E13").Value = '  -0.06%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.7170'
$ws.Range("E14").Value = '  +0.01%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '91.16'
$ws.Range("E15").Value = '  -0.39%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.319.18'
$ws.Range("E16").Value = '  -0.22%  '

# Row 17
$ws.Range("E17").Value = '  +1.88%  '

# Row 18
$ws.Range("E18").Value = '  +0.17%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '240.47'
$ws.Range("E19").Value = '  -1.77%  '

# Row 20
$ws.Range("E20").Value = '  +0.10%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '2.118.21'
$ws.Range("E21").Value = '  -0.42%  '

# Row 22
$ws.Range("E22").Value = '  +0.12%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.751'
$ws.Range("E23").Value = '  -2.33%  '

# Row 24
$ws.Range("E24").Value = '  +0.21%  '

# Row 25
$ws.Range("E25").Value = '  -2.60%  '

# Row 26
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.032'
$ws.Range("E26").Value = '  +0.04%  '

# Row 27
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '162.28'
$ws.Range("E27").Value = '  -0.83%  '

# Row 28
$ws.Range("E28").Value = '  -0.33%  '

# Row 29
$ws.Range("E29").Value = '  -0.17%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.406'
$ws.Range("E30").Value = '  -0.43%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.323'
$ws.Range("E31").Value = '  -0.03%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.281'
$ws.Range("E32").Value = '  -1.37%  '

# Row 33
$ws.Range("E33").Value = '  +3.03%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.944'
$ws.Range("E34").Value = '  +0.83%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7522'
$ws.Range("E35").Value = '  -2.51%  '

# Row 36
$ws.Range("E36").Value = '  +0.14%  '

# Row 37
$ws.Range("E37").Value = '  +0.08%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01885'
$ws.Range("E38").Value = '  +0.90%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.242.13'
$ws.Range("E39").Value = '  +6.56%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.734'
$ws.Range("E40").Value = '  +0.47%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.478'
$ws.Range("E41").Value = '  +1.35%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.8919'
$ws.Range("E42").Value = '  -0.07%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '72.31'
$ws.Range("E43").Value = '  -1.77%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '108.58'
$ws.Range("E44").Value = '  +4.35%  '

# Row 45
$ws.Range("E45").Value = '  +0.15%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.014.65'
$ws.Range("E46").Value = '  -0.47%  '

# Row 47
$ws.Range("E47").Value = '  +7.29%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.798'
$ws.Range("E48").Value = '  -0.37%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.5200'
$ws.Range("E49").Value = '  +0.10%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.451'
$ws.Range("E50").Value = '  +0.47%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.4335'
$ws.Range("E51").Value = '  +0.39%  '
